# Added the test-cases for MTTR.
# - New column entry "Parts" added in N3 (new shared string).
# - Sheet view scrolled/re-selected to reflect the newly added cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data cell: N3 = "Parts" (adds a new shared string, reuses the
# existing wrap-text cell style used throughout the sheet).
$ws.Range("N3").Value = "Parts"
$ws.Range("N3").WrapText = $true

# Update the view: scroll so column L is the left-most visible column,
# and select N3 as the active cell.
$win = $excel.ActiveWindow
$win.ScrollColumn = 12
$win.ScrollRow = 1
$ws.Range("N3").Select()
